$wb = $excel.ActiveWorkbook

# Physical worksheet slots (1st tab / 2nd tab), independent of their current names.
$slot1 = $wb.Worksheets.Item(1)   # currently "hotel_info"  -> becomes "review_info"
$slot2 = $wb.Worksheets.Item(2)   # currently "review_info" -> becomes "hotel_info"

# Wipe both sheets' existing contents - each slot gets entirely new data below.
$slot1.UsedRange.ClearContents()
$slot2.UsedRange.ClearContents()

# Swap the tab names (via a temporary name to avoid a collision).
$slot1.Name = "__tmp_swap__"
$slot2.Name = "hotel_info"
$slot1.Name = "review_info"

# --- slot1 (now "review_info"): write the review_info header row ---
$reviewHeaders = @(
    "STR","reviewer_ID","reviewer_name","Review_ID","Date_of_scraping","ReviewURL",
    "Tripadvisor_gcode","Tripadvisor_dcode","Tripadvisor_rcode","review_date","review_title",
    "review_content","review_rating","trip_month","trip_purpose","value","rooms","Location",
    "Cleanliness","Sleep Quality","Service","Picture(yes=1)","respondent","response_date","response_text"
)
for ($c = 1; $c -le $reviewHeaders.Count; $c++) {
    $slot1.Cells.Item(1, $c).Value = $reviewHeaders[$c - 1]
}

# --- slot2 (now "hotel_info"): write the hotel_info data, with a new "State" column
#     inserted right after "Hotel_Name" and before "City" ---
$hotelHeaders = @("STR","Hotel_Name","State","City","Zip","TA_ReviewURL","Tripadvisor_Hotel_Name","English_Reviews_num","Local_Rank","Total_Reviews_num")
for ($c = 1; $c -le $hotelHeaders.Count; $c++) {
    $slot2.Cells.Item(1, $c).Value = $hotelHeaders[$c - 1]
}

$slot2.Cells.Item(2, 1).Value = 40360
$slot2.Cells.Item(2, 2).Value = "Hotel St Helene"
$slot2.Cells.Item(2, 3).Value = "Louisiana"
$slot2.Cells.Item(2, 4).Value = "New Orleans"
$slot2.Cells.Item(2, 5).Value = 70130
$slot2.Cells.Item(2, 6).Value = "https://www.tripadvisor.com/Hotel_Review-g60864-d76995-Reviews-Hotel_St_Helene-New_Orleans_Louisiana.html"
$slot2.Cells.Item(2, 7).Value = "Hotel St. Helene"

# These three columns hold digit-strings in the source data ("374","119","383"), not
# numbers - force a Text format before assigning so they round-trip as shared strings.
$slot2.Cells.Item(2, 8).NumberFormat = "@"
$slot2.Cells.Item(2, 8).Value = "374"
$slot2.Cells.Item(2, 9).NumberFormat = "@"
$slot2.Cells.Item(2, 9).Value = "119"
$slot2.Cells.Item(2, 10).NumberFormat = "@"
$slot2.Cells.Item(2, 10).Value = "383"
